# Update the "想去人数" (attendance count) figures for a few events on the
# "展览" (exhibition) sheet and the "全部类型" (all-types) sheet.
#
# 展览 sheet (F column): row2 140->142, row3 454->455, row6 9->10, row9 109->124
# 全部类型 sheet (F column): row2 140->142, row4 454->455, row7 9->10, row10 109->124

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 142
$wsExhibit.Range("F3").Value = 455
$wsExhibit.Range("F6").Value = 10
$wsExhibit.Range("F9").Value = 124

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 142
$wsAll.Range("F4").Value = 455
$wsAll.Range("F7").Value = 10
$wsAll.Range("F10").Value = 124
